# Auto-update draw results: append the 2025-12-24 Pick 4 row to the
# results table (mirrors the existing A:E columns -> Date, Game, Phase,
# Result, InsertedAt) and keep everything stored as plain text, exactly
# like every other row already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# Leading apostrophes force text interpretation so date-looking /
# number-looking values ("2025-12-24", "251224") are stored as literal
# text instead of being reinterpreted as a date serial / number, just
# like the pre-existing rows.
$ws.Cells.Item($row, 1).Value = "'2025-12-24"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251224"
$ws.Cells.Item($row, 4).Value = "3-6-1-1"
$ws.Cells.Item($row, 5).Value = "2025-12-24T21:39:34.345+04:00"

# The table stores dates/phase-codes/results as text (not as real
# numbers/dates), so reset the new row back to the default "Normal"
# style/number format after the writes -- this clears any automatic
# number/date reinterpretation Excel applied while keeping the cell
# values as plain text, matching the rest of the sheet.
$ws.Range("A99:E99").Style = "Normal"
